$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 17.833334
$ws.Range("I8").Value = 17.833334
$ws.Range("K8").Value = 53.500002
$ws.Range("M8").Value = 85.49999800000001

# Row 62
$ws.Range("H62").Value = 6850.4707
$ws.Range("I62").Value = 5829.857
$ws.Range("J62").Value = 7564.9
$ws.Range("K62").Value = 5829.857
$ws.Range("L62").Value = 7564.9
$ws.Range("M62").Value = -5205.857
$ws.Range("N62").Value = -8812.9

# Row 65
$ws.Range("H65").Value = 6850.4707
$ws.Range("I65").Value = 5829.857
$ws.Range("J65").Value = 7564.9
$ws.Range("K65").Value = 29149.285
$ws.Range("L65").Value = 37824.5
$ws.Range("M65").Value = -26029.285
$ws.Range("N65").Value = -44064.5

# Row 138
$ws.Range("H138").Value = 7434.0713
$ws.Range("J138").Value = 7654.0347
$ws.Range("L138").Value = 22962.1041
$ws.Range("N138").Value = -33242.1041

# Row 141
$ws.Range("H141").Value = 606
$ws.Range("I141").Value = 583.375
$ws.Range("J141").Value = 666.3333
$ws.Range("K141").Value = 1750.125
$ws.Range("L141").Value = 1998.9999
$ws.Range("M141").Value = 3429.875
$ws.Range("N141").Value = -12358.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2635.7058
$ws.Range("I2").Value = 2021.9286
$ws.Range("K2").Value = 2021.9286
$ws.Range("M2").Value = -1908.9286

# Row 32
$ws.Range("H32").Value = 27408.295
$ws.Range("I32").Value = 19545.455
$ws.Range("J32").Value = 30029.242
$ws.Range("K32").Value = 19545.455
$ws.Range("L32").Value = 30029.242
$ws.Range("M32").Value = -19258.455
$ws.Range("N32").Value = -30603.242

# Row 45
$ws.Range("H45").Value = 2282.85
$ws.Range("I45").Value = 1471.5
$ws.Range("K45").Value = 1471.5
$ws.Range("M45").Value = -1094.5

# Row 74
$ws.Range("H74").Value = 2904.875
$ws.Range("I74").Value = 984.5
$ws.Range("K74").Value = 984.5
$ws.Range("M74").Value = -110.5

# Row 77
$ws.Range("H77").Value = 2904.875
$ws.Range("I77").Value = 984.5
$ws.Range("K77").Value = 4922.5
$ws.Range("M77").Value = -554.5

# Row 116
$ws.Range("H116").Value = 2635.7058
$ws.Range("I116").Value = 2021.9286
$ws.Range("K116").Value = 2021.9286
$ws.Range("M116").Value = 272.0714

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2635.7058
$ws.Range("I3").Value = 2021.9286
$ws.Range("K3").Value = 2021.9286
$ws.Range("M3").Value = -1907.9286

# Row 20
$ws.Range("H20").Value = 3077.1875
$ws.Range("I20").Value = 2777.7778
$ws.Range("K20").Value = 2777.7778
$ws.Range("M20").Value = -2530.7778

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 4709.846
$ws.Range("I58").Value = 1313.1428
$ws.Range("K58").Value = 1313.1428
$ws.Range("M58").Value = -1110.1428

# Row 94
$ws.Range("H94").Value = 982.7143
$ws.Range("I94").Value = 690
$ws.Range("K94").Value = 690
$ws.Range("M94").Value = -239

# Row 132
$ws.Range("H132").Value = 2929.3125
$ws.Range("I132").Value = 3032.8572
$ws.Range("J132").Value = 2731.6365
$ws.Range("K132").Value = 9098.571599999999
$ws.Range("L132").Value = 8194.9095
$ws.Range("M132").Value = -6568.571599999999
$ws.Range("N132").Value = -13254.9095

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()

# Row 136
$ws.Range("H136").Value = 4709.846
$ws.Range("I136").Value = 1313.1428
$ws.Range("K136").Value = 3939.4284
$ws.Range("M136").Value = -1389.4284

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1496.1794
$ws.Range("I131").Value = 1113.5714
$ws.Range("J131").Value = 1579.875
$ws.Range("K131").Value = 3340.7142
$ws.Range("L131").Value = 4739.625
$ws.Range("M131").Value = 1699.2858
$ws.Range("N131").Value = -14819.625

# Row 132
$ws.Range("H132").Value = 8499
$ws.Range("J132").Value = 11998.333
$ws.Range("L132").Value = 107984.997
$ws.Range("N132").Value = -113044.997

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 3872.5715
$ws.Range("I113").Value = 2343.4
$ws.Range("J113").Value = 4722.1113
$ws.Range("K113").Value = 2343.4
$ws.Range("L113").Value = 4722.1113
$ws.Range("M113").Value = -173.4000000000001
$ws.Range("N113").Value = -9062.1113

# Row 122
$ws.Range("H122").Value = 483974.4
$ws.Range("I122").Value = 77137.71000000001
$ws.Range("K122").Value = 231413.13
$ws.Range("M122").Value = -228963.13

# Row 132
$ws.Range("H132").Value = 3564.5
$ws.Range("I132").Value = 2636.9546
$ws.Range("K132").Value = 7910.8638
$ws.Range("M132").Value = -5380.8638

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3997
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 4995.5
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 4995.5
$ws.Range("M22").Value = -1705
$ws.Range("N22").Value = -5585.5

# Row 27
$ws.Range("H27").Value = 3997
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 4995.5
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 4995.5
$ws.Range("M27").Value = -1893
$ws.Range("N27").Value = -5209.5

# Row 82
$ws.Range("H82").Value = 1923.6111
$ws.Range("I82").Value = 2509.889
$ws.Range("J82").Value = 1337.3334
$ws.Range("K82").Value = 2509.889
$ws.Range("L82").Value = 1337.3334
$ws.Range("M82").Value = -2148.889
$ws.Range("N82").Value = -2059.3334

# Row 85
$ws.Range("H85").Value = 1923.6111
$ws.Range("I85").Value = 2509.889
$ws.Range("J85").Value = 1337.3334
$ws.Range("K85").Value = 2509.889
$ws.Range("L85").Value = 1337.3334
$ws.Range("M85").Value = -1261.889
$ws.Range("N85").Value = -3833.3334

# Row 132
$ws.Range("H132").Value = 3711.111
$ws.Range("I132").Value = 1936.3636
$ws.Range("K132").Value = 5809.0908
$ws.Range("M132").Value = -3279.0908

# Row 136
$ws.Range("H136").Value = 3100.8333
$ws.Range("I136").Value = 3100.8333
$ws.Range("K136").Value = 9302.499899999999
$ws.Range("M136").Value = -6752.499899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 4217.25

# Row 119
$ws.Range("H119").Value = 45000
$ws.Range("J119").Value = 45000
$ws.Range("L119").Value = 45000
$ws.Range("N119").Value = -54676

# Row 132
$ws.Range("H132").Value = 1981.9445
$ws.Range("I132").Value = 1206.1538
$ws.Range("K132").Value = 3618.4614
$ws.Range("M132").Value = -1088.4614
